# Rule edit: Mead Flood Control Release / Excel IC file edits - added TARV
#
# The IC files were updated to include a TARV (Target Annual Release
# Volume, MTOMRunType.TargetAnnualReleaseVolume_Input) for 2020, which was
# needed for simulations starting in August and September.
#
# This sets CoordinatedOps!E4 (the 2020-01-01 row) from the blank "NaN"
# placeholder to the numeric TARV value 9000000, and updates the saved
# cursor/selection state on the affected sheets to reflect where the user
# was working when the file was saved.

$wb = $excel.ActiveWorkbook

# --- CoordinatedOps sheet: set the TARV value for 2020 ---
$wsCoord = $wb.Worksheets.Item("CoordinatedOps")
$wsCoord.Activate() | Out-Null
$wsCoord.Range("E4").Value = 9000000
$wsCoord.Range("E5").Select() | Out-Null

# --- Reservoirs sheet: update the saved selection state ---
$wsRes = $wb.Worksheets.Item("Reservoirs")
$wsRes.Activate() | Out-Null
$wsRes.Range("S15:S26").Select() | Out-Null

# Restore CoordinatedOps as the active sheet (it was the tabSelected sheet)
$wsCoord.Activate() | Out-Null
